# Scheduled runner update for the crafting-leve "Profits" sheets
# (FFXIV gathering/crafting price tracker workbook).
#
# Summary of the change (see commit diff):
#   - ALC tab, row 138: refresh currentAveragePrice/currentAveragePriceHQ/
#     LevePriceHQ columns with new market-board data; LeveProfitHQ (N) no
#     longer computable so its cell is cleared entirely.
#   - BSM tab, row 99: prices dropped to 0 (no longer on market); profit
#     column (M) cleared entirely.
#   - BSM tab, rows 117-141: price-tracking columns (H:L, and for a couple
#     of rows M/N too) backfilled with freshly pulled data.
#   - CRP tab, rows 129-141 (except 136, already current): the stale
#     price-tracking columns (H:N) are removed outright so they will be
#     repopulated on the next run.

$wb = $excel.ActiveWorkbook

function Set-Cells {
    param($ws, $Row, $Values)
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
}

function Clear-Cells {
    param($ws, $Row, $Cols)
    foreach ($col in $Cols) {
        $ws.Range("$col$Row").ClearContents()
    }
}

# ---------------------------------------------------------------------
# ALC — row 138
# ---------------------------------------------------------------------
$wsALC = $wb.Worksheets.Item("ALC")

Set-Cells $wsALC 138 @{
    H = 1600
    J = 0
    L = 0
}
Clear-Cells $wsALC 138 @("N")

# ---------------------------------------------------------------------
# BSM — row 99
# ---------------------------------------------------------------------
$wsBSM = $wb.Worksheets.Item("BSM")

Set-Cells $wsBSM 99 @{
    H = 0
    I = 0
    K = 0
}
Clear-Cells $wsBSM 99 @("M")

# ---------------------------------------------------------------------
# BSM — rows 117-141: backfill currentAveragePrice* / LevePrice* /
# LeveProfit* columns that had been blank.
# ---------------------------------------------------------------------
$bsmZeroRows = @(117, 118, 119, 120, 122, 123, 124, 125, 126, 127, 128, `
    129, 130, 131, 132, 133, 135, 137, 138, 139, 141)

foreach ($row in $bsmZeroRows) {
    Set-Cells $wsBSM $row @{
        H = 0
        I = 0
        J = 0
        K = 0
        L = 0
    }
}

# Row 134 gets real (non-zero) figures, including profit columns.
Set-Cells $wsBSM 134 @{
    H = 1111.1111
    I = 833.3333
    J = 1666.6666
    K = 2499.9999
    L = 4999.9998
    M = 35.0001000000002
    N = -10069.9998
}

# Row 140 gets real figures too (no M value here).
Set-Cells $wsBSM 140 @{
    H = 75000
    I = 0
    J = 75000
    K = 0
    L = 75000
    N = -85360
}

# ---------------------------------------------------------------------
# CRP — rows 129-141 (except 136, which already reflects current data):
# drop the stale price-tracking columns so they get recomputed fresh.
# ---------------------------------------------------------------------
$wsCRP = $wb.Worksheets.Item("CRP")

$crpRowsHL = @(129, 130, 133, 137, 138, 139, 140, 141)
foreach ($row in $crpRowsHL) {
    Clear-Cells $wsCRP $row @("H", "I", "J", "K", "L")
}

Clear-Cells $wsCRP 131 @("H", "I", "J", "K", "L", "N")
Clear-Cells $wsCRP 132 @("H", "I", "J", "K", "L", "N")
Clear-Cells $wsCRP 134 @("H", "I", "J", "K", "L", "M")
Clear-Cells $wsCRP 135 @("H", "I", "J", "K", "L", "N")
# Row 136 intentionally left untouched.
